# Apply the LOM3066 workbook edit:
#  1. Update "Semestre ideal" value from "EM-8" to "EF-8,EM-8"
#  2. Remove the trailing "Requisitos:" block (rows 22-24)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Semestre ideal:" row (B9/C9) text.
$ws.Range("B9").Value = "EF-8,EM-8"
$ws.Range("C9").Value = "EF-8,EM-8"

# 2. Delete the last three rows (Requisitos: + the two requirement lines),
#    which also shrinks the used range from A1:C24 to A1:C21 and drops the
#    now-unreferenced shared strings.
$ws.Rows("22:24").Delete()
